$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits ---
# B2 holds a numeric-looking value but must stay text (matches source type),
# so force the cell to Text format before assigning - otherwise Excel
# auto-coerces "2" into the number 2.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("D2").ClearContents()
$ws.Range("G2").Value = 289.5
$ws.Range("H2").Value = 295
$ws.Range("I2").Value = 317
$ws.Range("J2").Value = 317
$ws.Range("K2").Value = -2.38
$ws.Range("L2").Value = 1190
$ws.Range("U2").Value = 5.441
$ws.Range("V2").Value = 0.08510871265446582
$ws.Range("W2").Value = -0.519039807206527
$ws.Range("X2").Value = 0.03554283097828655
$ws.Range("Y2").Value = -0.5545826381848136
$ws.Range("Z2").Value = 0.002560819462227913
$ws.Range("AA2").Value = -0.8107416879795395
$ws.Range("AB2").Value = 0.03575198178524969
$ws.Range("AC2").Value = -0.8464936697647891
$ws.Range("AD2").Value = 0.293
$ws.Range("AF2").Value = 0.293
$ws.Range("AG2").Value = -5.148
$ws.Range("AH2").Value = 0.004562228485122152
$ws.Range("AI2").Value = 0.04903765690376569
$ws.Range("AJ2").Value = -0.08757782994794325
$ws.Range("AK2").Value = -9.640449438202234
$ws.Range("AL2").Value = 0.028
$ws.Range("AM2").Value = 0.026
$ws.Range("AN2").Value = -0.5907258064516129
$ws.Range("AO2").Value = -22.64285714285714
$ws.Range("AP2").Value = 10.37903225806452
$ws.Range("AQ2").Value = -24.38461538461538

# --- Row 3 edits ---
$ws.Range("B3").Value = "Red Light Holland Corp. (CNSX:TRIP)"
$ws.Range("D3").ClearContents()
$ws.Range("G3").Value = -5.5
$ws.Range("H3").Value = -0
$ws.Range("I3").Value = -0
$ws.Range("J3").Value = -0
$ws.Range("K3").Value = -1.25
$ws.Range("L3").Value = 625
$ws.Range("U3").Value = 5
$ws.Range("V3").Value = 0.09090909090909091
$ws.Range("W3").Value = 0.6188118811881188
$ws.Range("X3").Value = 0.03533994470231133
$ws.Range("Y3").Value = 0.5834719364858074
$ws.Range("Z3").Value = 0.001706484641638225
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.03531845015025002
$ws.Range("AC3").Value = -0.03531845015025002
$ws.Range("AD3").Value = 0.089
$ws.Range("AF3").Value = 0.089
$ws.Range("AG3").Value = -4.911
$ws.Range("AH3").Value = 0.001615567536168745
$ws.Range("AI3").Value = 0.01679562181543687
$ws.Range("AJ3").Value = -0.09804547904729581
$ws.Range("AK3").Value = -16.42474916387958
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = -0.002
$ws.Range("AN3").ClearContents()
$ws.Range("AO3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AQ3").Value = -0

# --- Row 4: new row ---
$ws.Range("A4").Value = "Canada"
$ws.Range("B4").Value = "Harrys Manufacturing Inc. (CNSX:HARY)"
$ws.Range("C4").Value = "Tobacco"
$ws.Range("K4").Value = -1.13
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.441
$ws.Range("V4").Value = 0.04938409854423292
$ws.Range("W4").Value = -1.656891495601173
$ws.Range("X4").Value = 0.03574571725426177
$ws.Range("Y4").Value = -1.692637212855434
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = -1.621483375959079
$ws.Range("AB4").Value = 0.03618551342024936
$ws.Range("AC4").Value = -1.657668889379328
$ws.Range("AD4").Value = 0.204
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.204
$ws.Range("AG4").Value = -0.237
$ws.Range("AH4").Value = 0.02233413619443836
$ws.Range("AI4").Value = 0.3017751479289941
$ws.Range("AJ4").Value = -0.08757782994794325
$ws.Range("AK4").Value = -1.008510638297873
$ws.Range("AL4").Value = 0.028
$ws.Range("AM4").Value = 0.028
$ws.Range("AN4").Value = -0.4112903225806451
$ws.Range("AO4").Value = -22.64285714285714
$ws.Range("AP4").Value = 0.4778225806451613
$ws.Range("AQ4").Value = -22.64285714285714
